# Updated cryptos list -- refresh Price / Volume(1h) figures, and fix two
# rows whose coin order was swapped in the source feed (Filecoin/ImmutableX
# and Decentraland/EnergySwap).
#
# All D/E (and the occasional B/C) cells in this sheet are stored as plain
# text (inline strings) -- many of the "Price" values look numeric
# ("1.012", "21.00", ...) but must stay text so formatting like trailing
# zeros survives. Setting .Value directly would let Excel's COM layer
# auto-coerce number-looking strings into real numbers (and would also
# stamp a new NumberFormat style on the cell). Force the cell to Text
# first, write the value, then restore the (unstyled) "Normal" cell style
# so no stray formatting is introduced -- matches the original file, where
# none of these cells carry an explicit style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '27.856.47'
Set-TextValue 'E2' '  +2.78%  '
Set-TextValue 'D3' '1.868.51'
Set-TextValue 'E3' '  +1.05%  '
Set-TextValue 'D4' '1.012'
Set-TextValue 'E4' '  -0.41%  '
Set-TextValue 'D5' '312.94'
Set-TextValue 'E5' '  +1.20%  '
Set-TextValue 'E6' '  -0.38%  '
Set-TextValue 'D7' '0.4815'
Set-TextValue 'E7' '  +1.09%  '
Set-TextValue 'D8' '0.3819'
Set-TextValue 'E8' '  +3.64%  '
Set-TextValue 'D9' '0.07365'
Set-TextValue 'E9' '  +1.77%  '
Set-TextValue 'D10' '0.9375'
Set-TextValue 'E10' '  +0.66%  '
Set-TextValue 'D11' '21.00'
Set-TextValue 'E11' '  +5.76%  '
Set-TextValue 'D12' '0.07800'
Set-TextValue 'E12' '  +0.15%  '
Set-TextValue 'D13' '1.911.13'
Set-TextValue 'E13' '  +3.31%  '
Set-TextValue 'D14' '5.481'
Set-TextValue 'E14' '  +1.69%  '
Set-TextValue 'E15' '  +1.77%  '
Set-TextValue 'D16' '90.68'
Set-TextValue 'E16' '  +1.56%  '
Set-TextValue 'D17' '1.013'
Set-TextValue 'E17' '  -0.48%  '
Set-TextValue 'D18' '0.000008868'
Set-TextValue 'E18' '  +2.02%  '
Set-TextValue 'D19' '1.009'
Set-TextValue 'E19' '  -0.50%  '
Set-TextValue 'D20' '28.064.53'
Set-TextValue 'E20' '  +3.43%  '
Set-TextValue 'D21' '14.79'
Set-TextValue 'E21' '  +1.30%  '
Set-TextValue 'D22' '5.120'
Set-TextValue 'E22' '  +1.21%  '
Set-TextValue 'D23' '2.132.81'
Set-TextValue 'E23' '  +2.86%  '
Set-TextValue 'E24' '  +1.54%  '
Set-TextValue 'D25' '1.935'
Set-TextValue 'E25' '  -0.28%  '
Set-TextValue 'D26' '156.28'
Set-TextValue 'E26' '  +2.11%  '
Set-TextValue 'D27' '18.55'
Set-TextValue 'E27' '  +1.10%  '
Set-TextValue 'D28' '2.049'
Set-TextValue 'E28' '  +3.19%  '
Set-TextValue 'D29' '115.86'
Set-TextValue 'E29' '  +0.99%  '
Set-TextValue 'D30' '4.970'
Set-TextValue 'E30' '  +0.85%  '
Set-TextValue 'D31' '0.08908'
Set-TextValue 'E31' '  +0.40%  '
Set-TextValue 'D32' '3.329'
Set-TextValue 'E32' '  +0.94%  '
Set-TextValue 'D33' '1.218'
Set-TextValue 'E33' '  +3.06%  '

# Row 34/35: Filecoin and ImmutableX swap places (with new Price/Volume).
Set-TextValue 'B34' 'ImmutableX'
Set-TextValue 'C34' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D34' '0.7633'
Set-TextValue 'E34' '  +3.57%  '
Set-TextValue 'B35' 'Filecoin'
Set-TextValue 'C35' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D35' '4.658'
Set-TextValue 'E35' '  +3.01%  '

Set-TextValue 'D36' '2.726'
Set-TextValue 'E36' '  +1.97%  '
Set-TextValue 'D37' '1.134'
Set-TextValue 'E37' '  +1.67%  '
Set-TextValue 'D38' '0.02042'
Set-TextValue 'E38' '  +3.26%  '
Set-TextValue 'D39' '0.5632'
Set-TextValue 'E39' '  +6.72%  '
Set-TextValue 'D40' '0.05371'
Set-TextValue 'E40' '  +2.00%  '
Set-TextValue 'E41' '  +0.40%  '
Set-TextValue 'D42' '7.066'
Set-TextValue 'E42' '  +0.27%  '
Set-TextValue 'D43' '8.549'
Set-TextValue 'E43' '  +3.09%  '
Set-TextValue 'D44' '0.1531'
Set-TextValue 'E44' '  +0.50%  '

# Row 45/46: Decentraland and EnergySwap swap places (with new Price/Volume).
Set-TextValue 'B45' 'EnergySwap'
Set-TextValue 'C45' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D45' '10.78'
Set-TextValue 'E45' '  +1.54%  '
Set-TextValue 'B46' 'Decentraland'
Set-TextValue 'C46' 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue 'D46' '0.4891'
Set-TextValue 'E46' '  +3.26%  '

Set-TextValue 'D47' '1.009'
Set-TextValue 'D48' '104.80'
Set-TextValue 'E48' '  +2.76%  '
Set-TextValue 'D49' '1.674'
Set-TextValue 'E49' '  +3.53%  '
Set-TextValue 'D50' '67.61'
Set-TextValue 'D51' '0.06110'
